$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 128
$ws.Range("A128").NumberFormat = "@"
$ws.Range("A128").Value = "2024-05-14"
$ws.Range("A128").ClearFormats()
$ws.Range("B128").Value = "10:29:43"
$ws.Range("C128").Value = "No atornilla clips"
$ws.Range("D128").Value = "-"
$ws.Range("E128").Value = "-"
$ws.Range("F128").Value = "-"
$ws.Range("G128").Value = "-"

# Row 129
$ws.Range("A129").NumberFormat = "@"
$ws.Range("A129").Value = "2024-05-14"
$ws.Range("A129").ClearFormats()
$ws.Range("B129").Value = "10:33:57"
$ws.Range("C129").Value = "-"
$ws.Range("D129").Value = "Cámara no detecta Power CP"
$ws.Range("E129").Value = "-"
$ws.Range("F129").Value = "-"
$ws.Range("G129").Value = "-"
